# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates ---
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Range("C2").Value = 6
$rushing.Range("E2").Value = 11

$rushing.Range("C4").Value = 165
$rushing.Range("D4").Value = 81
$rushing.Range("F4").Value = 30

$rushing.Range("C5").Value = 22
$rushing.Range("D5").Value = 20
$rushing.Range("E5").Value = 12

# --- Receiving sheet updates ---
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 32
$receiving.Range("D2").Value = 27

$receiving.Range("C3").Value = 27
$receiving.Range("D3").Value = 23

$receiving.Range("C5").Value = 67
$receiving.Range("D5").Value = 51
$receiving.Range("E5").Value = 17
$receiving.Range("F5").Value = 10

$receiving.Range("C6").Value = 77
$receiving.Range("D6").Value = 54

$receiving.Range("C7").Value = 70
$receiving.Range("D7").Value = 47
$receiving.Range("E7").Value = 33
$receiving.Range("G7").Value = 11
$receiving.Range("H7").Value = 7

$receiving.Range("C11").Value = 42
$receiving.Range("D11").Value = 33

$receiving.Range("C12").Value = 11
$receiving.Range("D12").Value = 8
